$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.370.04"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "1.839.58"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "  +1.36%  "
$ws.Range("D5").Value = "'315.02"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("D7").Value = "'0.4745"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "'0.07464"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'0.8854"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "'20.49"
$ws.Range("D12").Value = "1.867.11"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("D13").Value = "'0.07378"
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").Value = "'5.486"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "'93.31"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "'1.015"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "'0.000008850"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "'1.014"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "'14.84"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "27.389.08"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "'5.350"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "2.068.04"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").Value = "'1.915"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "'152.46"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").Value = "'18.65"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").Value = "'2.168"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'5.254"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "'118.09"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "'0.08973"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'0.7600"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "'1.179"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "'4.556"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").Value = "'2.943"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").Value = "'1.014"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("D38").Value = "'0.05383"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").Value = "'0.01963"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'3.004"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").Value = "'7.299"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5359"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.391"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").Value = "'0.1665"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "'8.554"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").Value = "'0.4982"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'10.58"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").Value = "'1.015"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").Value = "'105.17"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "'1.682"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'0.06323"
$ws.Range("E51").Value = "  +0.34%  "
